$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2286585365853659
$ws.Range("C2").Value = 0.5274390243902439
$ws.Range("J2").Value = 0.006097560975609756
$ws.Range("P2").Value = 0.1310975609756098
$ws.Range("S2").Value = 0.1067073170731707
$ws.Range("B3").Value = 0.00558659217877095
$ws.Range("C3").Value = 0.0335195530726257
$ws.Range("J3").Value = 0.02793296089385475
$ws.Range("P3").Value = 0.7094972067039106
$ws.Range("S3").Value = 0.223463687150838
$ws.Range("J4").Value = 0.02272727272727273
$ws.Range("P4").Value = 0.6590909090909091
$ws.Range("S4").Value = 0.3181818181818182
$ws.Range("B6").Value = 0.07142857142857142
$ws.Range("D6").Value = 0.01260504201680672
$ws.Range("E6").Value = 0.004201680672268907
$ws.Range("F6").Value = 0.06302521008403361
$ws.Range("J6").Value = 0.2100840336134454
$ws.Range("O6").Value = 0.02941176470588235
$ws.Range("Q6").Value = 0.1470588235294118
$ws.Range("R6").Value = 0.07563025210084033
$ws.Range("S6").Value = 0.3865546218487395
$ws.Range("B7").Value = 0.125
$ws.Range("D7").Value = 0.0234375
$ws.Range("F7").Value = 0.046875
$ws.Range("J7").Value = 0.15625
$ws.Range("O7").Value = 0.01953125
$ws.Range("Q7").Value = 0.2109375
$ws.Range("R7").Value = 0.07421875
$ws.Range("S7").Value = 0.34375
$ws.Range("B8").Value = 0.1099476439790576
$ws.Range("D8").Value = 0.01570680628272251
$ws.Range("E8").Value = 0.002617801047120419
$ws.Range("F8").Value = 0.05759162303664921
$ws.Range("J8").Value = 0.08115183246073299
$ws.Range("O8").Value = 0.02879581151832461
$ws.Range("Q8").Value = 0.1884816753926702
$ws.Range("R8").Value = 0.06544502617801047
$ws.Range("S8").Value = 0.450261780104712
$ws.Range("B9").Value = 0.07798165137614679
$ws.Range("D9").Value = 0.009174311926605505
$ws.Range("F9").Value = 0.06422018348623854
$ws.Range("J9").Value = 0.1284403669724771
$ws.Range("O9").Value = 0.03669724770642202
$ws.Range("Q9").Value = 0.1926605504587156
$ws.Range("R9").Value = 0.06422018348623854
$ws.Range("S9").Value = 0.426605504587156
$ws.Range("B10").Value = 0.1075851393188854
$ws.Range("D10").Value = 0.02476780185758514
$ws.Range("F10").Value = 0.07198142414860681
$ws.Range("J10").Value = 0.1052631578947368
$ws.Range("O10").Value = 0.02786377708978328
$ws.Range("Q10").Value = 0.1934984520123839
$ws.Range("R10").Value = 0.07972136222910217
$ws.Range("S10").Value = 0.3893188854489164
$ws.Range("F11").Value = 0.002487562189054726
$ws.Range("G11").Value = 0.1467661691542289
$ws.Range("J11").Value = 0.09203980099502487
$ws.Range("K11").Value = 0.2238805970149254
$ws.Range("L11").Value = 0.5223880597014925
$ws.Range("S11").Value = 0.01243781094527363
$ws.Range("G12").Value = 0.7040358744394619
$ws.Range("J12").Value = 0.2331838565022422
$ws.Range("K12").Value = 0.01345291479820628
$ws.Range("L12").Value = 0.02242152466367713
$ws.Range("S12").Value = 0.02690582959641256
$ws.Range("G13").Value = 0.8070175438596491
$ws.Range("J13").Value = 0.1403508771929824
$ws.Range("S13").Value = 0.05263157894736842
$ws.Range("F15").Value = 0.04391891891891892
$ws.Range("H15").Value = 0.1385135135135135
$ws.Range("I15").Value = 0.1047297297297297
$ws.Range("J15").Value = 0.3006756756756757
$ws.Range("K15").Value = 0.09797297297297297
$ws.Range("M15").Value = 0.01351351351351351
$ws.Range("O15").Value = 0.09797297297297297
$ws.Range("S15").Value = 0.2027027027027027
$ws.Range("F16").Value = 0.02604166666666667
$ws.Range("I16").Value = 0.08854166666666667
$ws.Range("J16").Value = 0.4166666666666667
$ws.Range("K16").Value = 0.1510416666666667
$ws.Range("M16").Value = 0.01041666666666667
$ws.Range("O16").Value = 0.078125
$ws.Range("S16").Value = 0.1041666666666667
$ws.Range("F17").Value = 0.01345291479820628
$ws.Range("H17").Value = 0.1457399103139013
$ws.Range("I17").Value = 0.08295964125560538
$ws.Range("J17").Value = 0.4327354260089686
$ws.Range("K17").Value = 0.1143497757847534
$ws.Range("M17").Value = 0.02242152466367713
$ws.Range("O17").Value = 0.06278026905829596
$ws.Range("S17").Value = 0.1255605381165919
$ws.Range("F18").Value = 0.03888888888888889
$ws.Range("H18").Value = 0.15
$ws.Range("I18").Value = 0.1055555555555556
$ws.Range("J18").Value = 0.4055555555555556
$ws.Range("K18").Value = 0.1166666666666667
$ws.Range("M18").Value = 0.01111111111111111
$ws.Range("N18").Value = 0.005555555555555556
$ws.Range("O18").Value = 0.08888888888888889
$ws.Range("S18").Value = 0.07777777777777778
$ws.Range("F19").Value = 0.02084884586746091
$ws.Range("H19").Value = 0.1734921816827997
$ws.Range("I19").Value = 0.08637379002233805
$ws.Range("J19").Value = 0.3596425912137007
$ws.Range("K19").Value = 0.1295606850335071
$ws.Range("M19").Value = 0.02903946388682055
$ws.Range("N19").Value = 0.001489203276247208
$ws.Range("O19").Value = 0.08041697691734921
$ws.Range("S19").Value = 0.1191362620997766
